$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '247.51'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '1BNBBNB'
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.42'
$ws.Range('D3').ClearFormats()

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.638'
$ws.Range('D4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05609'
$ws.Range('D5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.395'
$ws.Range('D6').ClearFormats()

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.475'
$ws.Range('D7').ClearFormats()

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.074'
$ws.Range('D8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8019'
$ws.Range('D9').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07439'
$ws.Range('D11').ClearFormats()

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03191'
$ws.Range('D12').ClearFormats()

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.02964'
$ws.Range('D13').ClearFormats()

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09269'
$ws.Range('D14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001661'
$ws.Range('D15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.258'
$ws.Range('D16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04743'
$ws.Range('D17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.01179'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '17OneONEBestin24h'
$ws.Range('E18').ClearFormats()

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006262'
$ws.Range('D19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.001055'
$ws.Range('D20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.003814'
$ws.Range('D21').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0004778'
$ws.Range('D23').ClearFormats()

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.979'
$ws.Range('D24').ClearFormats()

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.109'
$ws.Range('D25').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1290'
$ws.Range('D27').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04193'
$ws.Range('D40').ClearFormats()

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('B41').ClearFormats()

$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('C41').ClearFormats()

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007011'
$ws.Range('D41').ClearFormats()

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.003501'
$ws.Range('D42').ClearFormats()

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('B43').ClearFormats()

$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('C43').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1047'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008852'
$ws.Range('D44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005684'
$ws.Range('D45').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6802'
$ws.Range('D47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.02948'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '47BOLOBOLOWorstin24h'
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('D49').ClearFormats()
